# Adding test cases for watch list
#
# - Adds a new test case row (TestCase_F14 / OPQA-1013) to the "Test Cases"
#   sheet, describing notification aggregation for post likes.
# - Marks the two preceding test cases (rows 13 & 14) as SKIP instead of PASS.
# - Highlights the header row with a yellow fill.
# - Updates the current selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Highlight header row (A1:E1) with yellow fill --------------------------
$ws.Range("A1:E1").Interior.Color = 65535

# --- Rows 13 & 14: Results column PASS -> SKIP ------------------------------
$ws.Range("E13").Value = "SKIP"
$ws.Range("E14").Value = "SKIP"

# --- New row 15: TestCase_F14 ------------------------------------------------
# Write the Description first, then TCID/Jira id/Runmode/Results, so that new
# shared-string entries are appended in the same order the workbook expects.
$ws.Range("C15").Value = "Verify that user is receiving notification when someone liked his post(aggregated notification)"
$ws.Range("A15").Value = "TestCase_F14"
$ws.Range("B15").Value = "OPQA-1013"
$ws.Range("D15").Value = "Y"
$ws.Range("E15").Value = "PASS"

# Match the formatting of the row above (borders / fills) for the new row.
$ws.Range("A14:E14").Copy() | Out-Null
$ws.Range("A15:E15").PasteSpecial(-4122) | Out-Null

# --- Update selection --------------------------------------------------------
$ws.Range("D7").Select() | Out-Null
